# Updates to Obsidian Search Engine Phase Overview Plans
# 1. Refresh the cached "datetimeFigureOut" footer-date field (30/08/2020 -> 09/10/2020)
#    on the Slide Master and on every Slide Layout.
# 2. Re-word the "Causal Regressions" label (shape "TextBox 11" on slide 1) to
#    "Causal Methods", written as two runs ("Causal " + "Methods") to mirror the
#    authored edit.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "09/10/2020"
        }
    }
}

# --- Slide Master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Slide Layouts ---
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# --- Slide 1: "Causal Regressions" -> "Causal " + "Methods" ---
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Causal Regressions") {
        $tr = $shp.TextFrame.TextRange
        # Retype just the "Regressions" portion (chars 8-19) as "Methods", leaving
        # "Causal " (chars 1-7) as-is; this naturally yields two runs: "Causal " + "Methods".
        $tr.Characters(8, 12).Text = "Methods"
    }
}
